$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input values (weight / recoil figures got lighter) ---
$ws.Range("D3").Value = 0.28999999999999998
$ws.Range("D4").Value = 0.28999999999999998

$ws.Range("D5").Value = 0.22
$ws.Range("E5").Value = -8
$ws.Range("F5").Value = -11

$ws.Range("D6").Value = 0.22
$ws.Range("E6").Value = -8
$ws.Range("F6").Value = -11

# --- Update hardcoded irl-weight (oz) figures ---
$ws.Range("P3").Value = 14.356517999999999
$ws.Range("P4").Value = 14.356517999999999
$ws.Range("P5").Value = 10.987850999999999
$ws.Range("P6").Value = 10.987850999999999

# --- Row 55: formulas now reference rows 4/7/9 directly instead of row 53 ---
foreach ($col in 4..16) {
    $letter = [char](64 + $col)
    $ws.Cells.Item(55, $col).Formula = "=" + $letter + "4+" + $letter + "7+" + $letter + "9"
}
$ws.Range("C55").Formula = "=C4+C7+C9"

# --- Row 56: formulas now reference rows 5/7/9 directly instead of row 54 ---
foreach ($col in 4..16) {
    $letter = [char](64 + $col)
    $ws.Cells.Item(56, $col).Formula = "=" + $letter + "5+" + $letter + "7+" + $letter + "9"
}
$ws.Range("C56").Formula = "=C5+C7+C9"

# --- Restore sheet view scroll/selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("S44").Select()
